$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.173.03"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "3.107.24"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "580.55"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "173.16"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.101.13"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "37.28"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "3.621.06"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "67.135.51"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "3.105.34"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "16.24"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "485.94"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "84.22"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "13.27"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "8.03"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "48.00"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "2.14"
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("D40").Value = "50.24"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Value = "8.68"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").Value = "2.854.25"
$ws.Range("E45").Value = "  +4.80%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "385.15"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "135.97"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "25.05"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  -0.50%  "
